$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-9 to the new serial value
$ws.Range("C2:C9").Value = 46074

# Rows 4-9 hold a rotating set of case records (A: Beteckning, B: Datum,
# F: Markägare, G: Area (ha)). The underlying source data shifted, so we
# rewrite each row's values to match the new snapshot.

# Row 4
$ws.Range("A4").Value = "A 6983-2023"
$ws.Range("B4").Value = 44967.68585648148
$ws.Range("F4").ClearContents()
$ws.Range("G4").Value = 5.4

# Row 5
$ws.Range("A5").Value = "A 25251-2025"
$ws.Range("B5").Value = 45800.50082175926
$ws.Range("F5").Value = "Kommuner"
$ws.Range("G5").Value = 0.7

# Row 6
$ws.Range("A6").Value = "A 35734-2023"
$ws.Range("B6").Value = 45147.89258101852
$ws.Range("F6").ClearContents()
$ws.Range("G6").Value = 5.9

# Row 7
$ws.Range("A7").Value = "A 5402-2026"
$ws.Range("B7").Value = 46050.49721064815
$ws.Range("F7").Value = "Kommuner"
$ws.Range("G7").Value = 0.7

# Row 8
$ws.Range("A8").Value = "A 25254-2025"
$ws.Range("B8").Value = 45800.50479166667
$ws.Range("G8").Value = 0.2

# Row 9
$ws.Range("A9").Value = "A 26074-2025"
$ws.Range("B9").Value = 45805.32366898148
$ws.Range("G9").Value = 1.3
